# Update column F (dSF) values to match repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    4  = 2
    5  = 1
    6  = -2
    7  = 2
    8  = -3
    9  = 5
    10 = 2
    14 = -5
    15 = 2
    16 = 3
    17 = 6
    19 = 7
    20 = -2
    21 = -3
    23 = -6
    26 = 2
    29 = -15
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
